$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Add the new "TROY" worksheet as the last tab
# ------------------------------------------------------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $last)
$ws.Name = "TROY"

# ------------------------------------------------------------------
# Cell values - entered in the same order the original author typed
# them so the shared-string table comes out in the same sequence.
# ------------------------------------------------------------------
$ws.Range("C6").Value = "Number of tasks"
$ws.Range("E6").Value = "Number of threads per task"
$ws.Range("F6").Value = "Machine"
$ws.Range("F7").Value = "painter"
$ws.Range("G6").Value = "Type"
$ws.Range("G7").Value = "matching"
$ws.Range("H6").Value = "TTC"
$ws.Range("B6").Value = "Backend"
$ws.Range("B7").Value = "BJ-SAGA"
$ws.Range("B8").Value = "BJ-Diane"
$ws.Range("B9").Value = "BJ-SAGA, BJ-Diane"
$ws.Range("D6").Value = "Number of cores"

$ws.Range("C7").Value = 8
$ws.Range("D7").Value = 16
$ws.Range("E7").Value = 2
$ws.Range("H7").Value = 529

$ws.Range("C8").Value = 8
$ws.Range("D8").Value = 16
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = "painter"
$ws.Range("G8").Value = "matching"
$ws.Range("H8").Value = 920

$ws.Range("C9").Value = 8
$ws.Range("D9").Value = 16
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = "painter"
$ws.Range("G9").Value = "matching"

# ------------------------------------------------------------------
# Column widths
# ------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 18.83203125
$ws.Columns.Item(5).ColumnWidth = 12.33203125

# ------------------------------------------------------------------
# Header row (row 6) formatting - built-in "Heading 1" cell style
# ------------------------------------------------------------------
$ws.Range("B6:H6").Style = "Heading 1"
$ws.Range("B6:H6").WrapText = $true
$ws.Rows.Item(6).RowHeight = 58

# ------------------------------------------------------------------
# Selection / active-cell bookkeeping for the new sheet
# ------------------------------------------------------------------
$ws.Range("H11").Select()

# Old "File Transfer" sheet loses the tab-selected flag and its
# selection moves
$ws1 = $wb.Worksheets.Item("File Transfer")
$ws1.Range("F16").Select()

# TROY becomes the active (selected) tab
$ws.Activate()
